$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by one column (N) that mirrors column M's formatting,
# adding the 2022 data point to each of the existing rows (3-7).

# Row 3: thin bottom-border separator row, style only (no value), same as M3.
$ws.Range("M3").Copy($ws.Range("N3"))

# Row 4: header row with year values.
$ws.Range("M4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = 2022

# Row 5: 2G coverage.
$ws.Range("M5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = 98.8

# Row 6: 3G coverage.
$ws.Range("M6").Copy($ws.Range("N6"))
$ws.Range("N6").Value = 98

# Row 7: 4G coverage.
$ws.Range("M7").Copy($ws.Range("N7"))
$ws.Range("N7").Value = 96.9

# Match the author's active selection at the time of saving.
$ws.Range("O4").Select() | Out-Null
